$d = $word.ActiveDocument

# 1) The publication line "Book chapter for Royal Society of Chemistry. " /
#    "In press" / "." was split across three separate runs. Collapse them
#    into a single run reading "Book chapter for Royal Society of Chemistry.
#    In press." (Find/Replace across the run boundary merges the runs,
#    inheriting the formatting of the first matched run.)
$d.Content.Find.Execute("Book chapter for Royal Society of Chemistry. In press.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Book chapter for Royal Society of Chemistry. In press.", 2)

# 2) Drop the stray leading-space run that preceded the "Presentations"
#    section heading run, leaving just the "Presentations" run.
$d.Content.Find.Execute(" Presentations", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Presentations", 2)
